# Update financial DB download file
# Sheet "2021": insert 4 new ANSP rows (IAA, LGS, NAVIAIR, Oro Navigacija)
# into the existing data block, shifting the rows for LFV and Skyguide
# down to their new positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2021")

# Row 4 currently holds LFV, row 5 holds Skyguide.
# Insert a new row at 4 so LFV moves to row 5, Skyguide moves to row 6.
$ws.Rows.Item(4).Insert()

# Insert three more rows at 6 (right after the now-shifted LFV row 5),
# pushing Skyguide from row 6 down to row 9.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

$ws.Cells.Item(4, 1).Value = 2021
$ws.Range("B4").Value = "IAA"
$ws.Range("C4").Value = "€"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Million"
$ws.Range("F4").Value = 162.256
$ws.Range("G4").Value = -144.513
$ws.Range("H4").Value = 17.743
$ws.Range("I4").Value = -10.673
$ws.Range("J4").Value = 7.07
$ws.Range("K4").Value = -1.378
$ws.Range("L4").Value = 11.379
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 141.487
$ws.Range("O4").Value = 275.802
$ws.Range("P4").Value = 417.289
$ws.Range("Q4").Value = 260.405
$ws.Range("R4").Value = 39.442
$ws.Range("S4").Value = 117.442
$ws.Range("T4").Value = 156.884
$ws.Range("U4").Value = 116.669
$ws.Range("V4").Value = 294.673731774996
$ws.Range("W4").Value = -23.654
$ws.Range("X4").Value = 70.04
$ws.Range("Y4").Value = -0.898
$ws.Range("Z4").Value = 45.488
$ws.Range("AA4").Value = -11.775
$ws.Range("AB4").Value = -35.429

$ws.Cells.Item(6, 1).Value = 2021
$ws.Range("B6").Value = "LGS"
$ws.Range("C6").Value = "€"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "Million"
$ws.Range("F6").Value = 19.797503
$ws.Range("G6").Value = -19.003318
$ws.Range("H6").Value = 0.794185
$ws.Range("I6").Value = -3.906371
$ws.Range("J6").Value = -3.112186
$ws.Range("K6").Value = -0.04245
$ws.Range("L6").Value = -3.154636
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 31.848425
$ws.Range("O6").Value = 7.521389
$ws.Range("P6").Value = 39.369814
$ws.Range("Q6").Value = 32.356607
$ws.Range("R6").Value = 4.322769
$ws.Range("S6").Value = 2.690438
$ws.Range("T6").Value = 7.013207
$ws.Range("U6").Value = 2.871421
$ws.Range("V6").Value = 55.1518774247739
$ws.Range("W6").Value = -2.262872
$ws.Range("X6").Value = -4.282818
$ws.Range("Y6").Value = -0.223366
$ws.Range("Z6").Value = -6.769056
$ws.Range("AA6").Value = -4.285703
$ws.Range("AB6").Value = -6.548575

$ws.Cells.Item(7, 1).Value = 2021
$ws.Range("B7").Value = "NAVIAIR"
$ws.Range("C7").Value = "DKK"
$ws.Range("D7").Value = 7.437
$ws.Range("E7").Value = "Million"
$ws.Range("F7").Value = 884.256
$ws.Range("G7").Value = -797.606
$ws.Range("H7").Value = 86.65
$ws.Range("I7").Value = -107.583
$ws.Range("J7").Value = 5.412
$ws.Range("K7").Value = -10.655
$ws.Range("L7").Value = -3.314
$ws.Range("M7").Value = 26.345
$ws.Range("N7").Value = 1298.281
$ws.Range("O7").Value = 1018.242
$ws.Range("P7").Value = 2316.523
$ws.Range("Q7").Value = 1082.781
$ws.Range("R7").Value = 337.249
$ws.Range("S7").Value = 896.493
$ws.Range("T7").Value = 1233.742
$ws.Range("U7").Value = 39.813
$ws.Range("V7").Value = 18.2192022126213
$ws.Range("W7").Value = -286.908
$ws.Range("X7").Value = -82.449
$ws.Range("Y7").Value = 398.75
$ws.Range("Z7").Value = 29.393
$ws.Range("AA7").Value = -82.449
$ws.Range("AB7").Value = -369.357

$ws.Cells.Item(8, 1).Value = 2021
$ws.Range("B8").Value = "Oro Navigacija"
$ws.Range("C8").Value = "€"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "Million"
$ws.Range("F8").Value = 26.898
$ws.Range("G8").Value = -18.108
$ws.Range("H8").Value = 8.79
$ws.Range("I8").Value = -3.906
$ws.Range("J8").Value = 4.884
$ws.Range("K8").Value = -0.014
$ws.Range("L8").Value = 4.87
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 45.145
$ws.Range("O8").Value = 14.703
$ws.Range("P8").Value = 59.848
$ws.Range("Q8").Value = 47.112
$ws.Range("R8").Value = 5.268
$ws.Range("S8").Value = 7.468
$ws.Range("T8").Value = 12.736
$ws.Range("U8").Value = 6.914
$ws.Range("V8").Value = 139.364369339518
$ws.Range("W8").Value = -0.491
$ws.Range("X8").Value = -3.239
$ws.Range("Y8").Value = -1.431
$ws.Range("Z8").Value = -5.161
$ws.Range("AA8").Value = -2.262
$ws.Range("AB8").Value = -2.753

